$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D12: TMAP Auto description - replace HTML <a> link with plain <br>-formatted URL
$ws.Range("D12").Value = "TMAP Auto에 대해 문의하는 경우 아래 페이지를 참고합니다`n<br>`n<br>https://tmaphelp.zendesk.com/hc/ko/sections/25814297789083-%EB%B3%BC%EB%B3%B4"

# D13: remove stray leading spaces before <br> on continuation lines
$ws.Range("D13").Value = "고객이 차량을 탑승할 때마다 중앙 화면이 초기화되는 경우, 차량을 타인에게 빌려주는 경우 등 운전자 모드 변경 방법에 대해 문의하는 경우 아래 내용을 참고합니다.`n<br>`n<br>1. 차량의 상단에서 밑으로 내려 알림 센터(Notification Center)를 내립니다.`n<br>2. 상단의 프로필을 누른 후 프로필을 변경합니다."

# A14: "티맵" -> "TMAP" in title
$ws.Range("A14").Value = "TMAP 운전 점수 문의를 주는 경우"
# D14: add missing <br> prefix to last bullet line
$ws.Range("D14").Value = "고객이 차량의 티맵 오토에 로그인하기 전 주행 기록을 로그인 후 주행 기록과 합치고자 하는 경우 아래 내용을 안내합니다.`n <br>`n<br>- 운전 점수는 차량 정보와 무관하게 티맵 계정 기반으로 관리되기에 두 데이터를 합칠 수 없는 점에 대해 안내 합니다.`n<br>- 자세한 내용은 티맵에 문의하도록 고객에게 안내합니다."

# A15: "플로앱" -> "Flo앱" in title
$ws.Range("A15").Value = "Flo앱이 동작하지 않는 경우"
# D15: "볼보" -> "Volvo"
$ws.Range("D15").Value = "고객이 예기치 않게 플로 앱이 재생되지 않는 경우 아래 사항을 확인합니다.`n<br>* 1분 미리 듣기만 되는 경우`n<br>`n<br>1. Volvo에서 차량 출고시 제공하는 이용권은 1년임을 안내 합니다.`n<br>2. 1년 이후 이용권 구입 시 이용 가능함을 안내 합니다.`n<br>3. 이용권 구입했음에도 동일한 경우 플로 앱에서 이용권 식별 후 없는 경우 플로 고객센터 안내합니다."

# Restore view state: scroll position and active cell selection
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("D9").Select()
